$p = $ppt.ActivePresentation

# --- Update the cached "datetimeFigureOut" date placeholder text on the
# slide master and every slide layout (mirrors Insert > Header & Footer >
# Date and time > Fixed "10/01/2023" -> "2023-10-05", applied to all). ---
$oldDate = "10/01/2023"
$newDate = "2023-10-05"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.HasTextFrame) {
                if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                    $shp.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

Update-DatePlaceholder $p.SlideMaster.Shapes

for ($j = 1; $j -le $p.SlideMaster.CustomLayouts.Count; $j++) {
    $lay = $p.SlideMaster.CustomLayouts.Item($j)
    Update-DatePlaceholder $lay.Shapes
}

# --- Reposition the "Oval 20" shape on slide 1 ---
# PowerPoint's Shape.Left/.Top are expressed in points (1 pt = 12700 EMU),
# so the target EMU offsets from the OOXML are converted to points. A
# small self-correcting nudge loop guards against float round-trip jitter
# landing one EMU off the intended value.
function Set-ShapeEmuPosition($shp, $targetXEmu, $targetYEmu) {
    $ptX = $targetXEmu / 12700.0
    $ptY = $targetYEmu / 12700.0
    $shp.Left = $ptX
    $shp.Top = $ptY

    $guardX = 0
    while (([Math]::Round($shp.Left * 12700)) -ne $targetXEmu -and $guardX -lt 20) {
        if (([Math]::Round($shp.Left * 12700)) -lt $targetXEmu) {
            $ptX = $ptX + 0.000005
        } else {
            $ptX = $ptX - 0.000005
        }
        $shp.Left = $ptX
        $guardX = $guardX + 1
    }

    $guardY = 0
    while (([Math]::Round($shp.Top * 12700)) -ne $targetYEmu -and $guardY -lt 20) {
        if (([Math]::Round($shp.Top * 12700)) -lt $targetYEmu) {
            $ptY = $ptY + 0.000005
        } else {
            $ptY = $ptY - 0.000005
        }
        $shp.Top = $ptY
        $guardY = $guardY + 1
    }
}

$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "Oval 20") {
        Set-ShapeEmuPosition $shp 7929943 435694
    }
}
